$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{Cell="AA2"; Value=7.5}
  @{Cell="G2"; Value=2.1}
  @{Cell="H2"; Value=3.6}
  @{Cell="I2"; Value=3.3}
  @{Cell="V2"; Value=9.5}
  @{Cell="X2"; Value=17}
  @{Cell="Z2"; Value=13}
  @{Cell="AF4"; Value=41}
  @{Cell="AG4"; Value=23}
  @{Cell="AH4"; Value=101}
  @{Cell="G4"; Value=1.27}
  @{Cell="I4"; Value=8.5}
  @{Cell="R4"; Value=1.8}
  @{Cell="S4"; Value=1.95}
  @{Cell="K5"; Value=8}
  @{Cell="AB8"; Value=15}
  @{Cell="AF8"; Value=21}
  @{Cell="G8"; Value=1.9}
  @{Cell="I8"; Value=3.75}
  @{Cell="K8"; Value=12}
  @{Cell="R8"; Value=1.75}
  @{Cell="S8"; Value=2}
  @{Cell="AA9"; Value=8.5}
  @{Cell="G9"; Value=1.65}
  @{Cell="P9"; Value=1.25}
  @{Cell="Q9"; Value=3.75}
  @{Cell="T9"; Value=10}
  @{Cell="W9"; Value=13}
  @{Cell="Y9"; Value=19}
  @{Cell="AB10"; Value=21}
  @{Cell="AD10"; Value=251}
  @{Cell="AF10"; Value=51}
  @{Cell="AG10"; Value=29}
  @{Cell="AI10"; Value=67}
  @{Cell="AJ10"; Value=51}
  @{Cell="G10"; Value=1.22}
  @{Cell="H10"; Value=6.25}
  @{Cell="I10"; Value=12}
  @{Cell="R10"; Value=1.91}
  @{Cell="S10"; Value=1.91}
  @{Cell="T10"; Value=10}
  @{Cell="U10"; Value=7.5}
  @{Cell="W10"; Value=8}
  @{Cell="Y10"; Value=23}
  @{Cell="Z10"; Value=21}
  @{Cell="AB13"; Value=17}
  @{Cell="AC13"; Value=100}
  @{Cell="AE13"; Value=7.4}
  @{Cell="AF13"; Value=14}
  @{Cell="AG13"; Value=11}
  @{Cell="AI13"; Value=30}
  @{Cell="AJ13"; Value=45}
  @{Cell="L13"; Value=1.44}
  @{Cell="M13"; Value=2.42}
  @{Cell="N13"; Value=2.27}
  @{Cell="O13"; Value=1.5}
  @{Cell="P13"; Value=1.52}
  @{Cell="Q13"; Value=2.2}
  @{Cell="R13"; Value=1.93}
  @{Cell="S13"; Value=1.7}
  @{Cell="T13"; Value=6.4}
  @{Cell="Z13"; Value=7}
  @{Cell="AA15"; Value=5.8}
  @{Cell="AB15"; Value=17}
  @{Cell="AC15"; Value=110}
  @{Cell="AD15"; Value=800}
  @{Cell="AE15"; Value=6.7}
  @{Cell="AG15"; Value=10.25}
  @{Cell="AH15"; Value=29}
  @{Cell="AI15"; Value=26}
  @{Cell="AJ15"; Value=45}
  @{Cell="G15"; Value=2.77}
  @{Cell="H15"; Value=2.92}
  @{Cell="I15"; Value=2.6}
  @{Cell="J15"; Value=1.09}
  @{Cell="K15"; Value=6}
  @{Cell="L15"; Value=1.47}
  @{Cell="M15"; Value=2.35}
  @{Cell="N15"; Value=2.32}
  @{Cell="O15"; Value=1.47}
  @{Cell="P15"; Value=1.52}
  @{Cell="Q15"; Value=2.22}
  @{Cell="R15"; Value=1.98}
  @{Cell="S15"; Value=1.65}
  @{Cell="T15"; Value=6.7}
  @{Cell="U15"; Value=12.5}
  @{Cell="V15"; Value=10.75}
  @{Cell="X15"; Value=29}
  @{Cell="Y15"; Value=50}
  @{Cell="Z15"; Value=6.6}
  @{Cell="AB16"; Value=13}
  @{Cell="AE16"; Value=8.75}
  @{Cell="AF16"; Value=17}
  @{Cell="AG16"; Value=11}
  @{Cell="AH16"; Value=45}
  @{Cell="AI16"; Value=30}
  @{Cell="G16"; Value=2.32}
  @{Cell="I16"; Value=3.25}
  @{Cell="Q16"; Value=2.7}
  @{Cell="R16"; Value=1.75}
  @{Cell="S16"; Value=1.95}
  @{Cell="T16"; Value=7.4}
  @{Cell="U16"; Value=11.5}
  @{Cell="W16"; Value=25}
  @{Cell="X16"; Value=19}
  @{Cell="Y16"; Value=28}
  @{Cell="AD22"; Value=501}
  @{Cell="AE22"; Value=10}
  @{Cell="AH22"; Value=51}
  @{Cell="AJ22"; Value=51}
  @{Cell="G22"; Value=1.85}
  @{Cell="H22"; Value=3.2}
  @{Cell="I22"; Value=4.75}
  @{Cell="J22"; Value=1.08}
  @{Cell="K22"; Value=8}
  @{Cell="N22"; Value=2.3}
  @{Cell="O22"; Value=1.6}
  @{Cell="U22"; Value=8}
  @{Cell="W22"; Value=15}
  @{Cell="X22"; Value=17}
  @{Cell="AF23"; Value=21}
  @{Cell="AG23"; Value=15}
  @{Cell="G23"; Value=1.95}
  @{Cell="I23"; Value=4.33}
  @{Cell="J23"; Value=1.07}
  @{Cell="K23"; Value=8.5}
  @{Cell="N23"; Value=2.08}
  @{Cell="O23"; Value=1.73}
  @{Cell="Z23"; Value=8.5}
  @{Cell="AD24"; Value=251}
  @{Cell="AE24"; Value=13}
  @{Cell="H24"; Value=3.2}
  @{Cell="N24"; Value=2.05}
  @{Cell="O24"; Value=1.75}
  @{Cell="AA25"; Value=6.5}
  @{Cell="AC25"; Value=41}
  @{Cell="AD25"; Value=201}
  @{Cell="AE25"; Value=8.5}
  @{Cell="AJ25"; Value=26}
  @{Cell="H25"; Value=3.25}
  @{Cell="L25"; Value=1.25}
  @{Cell="M25"; Value=3.75}
  @{Cell="N25"; Value=1.93}
  @{Cell="O25"; Value=1.93}
  @{Cell="P25"; Value=1.36}
  @{Cell="Q25"; Value=3}
  @{Cell="R25"; Value=1.7}
  @{Cell="S25"; Value=2.05}
  @{Cell="T25"; Value=10}
  @{Cell="X25"; Value=23}
  @{Cell="Y25"; Value=29}
  @{Cell="Z25"; Value=11}
  @{Cell="AA26"; Value=6.5}
  @{Cell="AF26"; Value=26}
  @{Cell="AI26"; Value=41}
  @{Cell="G26"; Value=1.8}
  @{Cell="I26"; Value=5}
  @{Cell="T26"; Value=7.5}
  @{Cell="U26"; Value=8.5}
  @{Cell="AA30"; Value=5.4}
  @{Cell="AB30"; Value=14.5}
  @{Cell="AC30"; Value=80}
  @{Cell="AD30"; Value=500}
  @{Cell="AE30"; Value=7.8}
  @{Cell="AF30"; Value=16.5}
  @{Cell="AG30"; Value=11.75}
  @{Cell="AH30"; Value=50}
  @{Cell="AI30"; Value=35}
  @{Cell="AJ30"; Value=45}
  @{Cell="G30"; Value=1.85}
  @{Cell="H30"; Value=3.15}
  @{Cell="I30"; Value=4}
  @{Cell="L30"; Value=1.38}
  @{Cell="M30"; Value=2.85}
  @{Cell="N30"; Value=2.2}
  @{Cell="O30"; Value=1.52}
  @{Cell="P30"; Value=1.45}
  @{Cell="Q30"; Value=2.27}
  @{Cell="R30"; Value=2.04}
  @{Cell="S30"; Value=1.7}
  @{Cell="T30"; Value=4.85}
  @{Cell="U30"; Value=6.5}
  @{Cell="V30"; Value=7.3}
  @{Cell="W30"; Value=12}
  @{Cell="X30"; Value=14}
  @{Cell="Y30"; Value=28}
  @{Cell="Z30"; Value=7.1}
  @{Cell="K33"; Value=17}
  @{Cell="AD35"; Value=600}
  @{Cell="AB36"; Value=10}
  @{Cell="AC36"; Value=23}
  @{Cell="AE36"; Value=17}
  @{Cell="AF36"; Value=19}
  @{Cell="G36"; Value=2.5}
  @{Cell="H36"; Value=3.8}
  @{Cell="N36"; Value=1.36}
  @{Cell="O36"; Value=3}
  @{Cell="P36"; Value=1.2}
  @{Cell="Q36"; Value=4.33}
  @{Cell="Y36"; Value=19}
  @{Cell="Z36"; Value=23}
)

foreach ($u in $updates) {
  $ws.Range($u.Cell).Value = $u.Value
}
